$wb = $excel.ActiveWorkbook

# --- Add the new "Invalid_Login" worksheet as the last tab ---
$ws2 = $wb.Worksheets.Item("Valid_Login")
$ws3 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "Invalid_Login"

# --- Row 1: header ---
$ws3.Range("A1").Value = "LoginX"
$ws3.Range("B1").Value = "Password"
$ws3.Range("C1").Value = "Type"

# --- Row 2: valid-looking login, wrong password style test ---
$ws3.Range("A2").Value = "LoginX"
$ws3.Range("B2").Value = "Pwd@abcd1"
$ws3.Range("C2").Value = "invalidLoginName"
$ws3.Hyperlinks.Add($ws3.Range("B2"), "mailto:Pwd@abcd1") | Out-Null
$ws3.Range("B2").Style = "Hyperlink"

# --- Row 3: null login name ---
$ws3.Range("B3").Value = "Pwd@abcd1"
$ws3.Range("C3").Value = "NullLoginName"
$ws3.Hyperlinks.Add($ws3.Range("B3"), "mailto:Pwd@abcd1") | Out-Null
$ws3.Range("B3").Style = "Hyperlink"

# --- Row 4: null password ---
$ws3.Range("A4").Value = "TEST"
$ws3.Range("C4").Value = "NullPwd"

# --- Blank-but-quote-prefixed cells, entered last (as in the source data) ---
$ws3.Range("A3").Value = "'"
$ws3.Range("B4").Value = "'"

# --- Selection / active-tab bookkeeping to match the final UI state ---
$ws2.Range("A1:D2").Select()
$ws3.Columns("C").Select()
$ws3.Activate()
